# Updated cryptos list on Tue Dec  5 18:47:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue 'D2' '43.953.69'
Set-TextValue 'E2' '  +5.34%  '
Set-TextValue 'D3' '2.297.09'
Set-TextValue 'E3' '  +3.31%  '
Set-TextValue 'E4' '  +0.24%  '
Set-TextValue 'D5' '232.41'
Set-TextValue 'E5' '  +0.47%  '
Set-TextValue 'D6' '0.628'
Set-TextValue 'E6' '  +0.61%  '
Set-TextValue 'D7' '61.93'
Set-TextValue 'E7' '  +2.05%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'D9' '0.419'
Set-TextValue 'E9' '  +4.45%  '
Set-TextValue 'D10' '0.0921'
Set-TextValue 'E10' '  +3.65%  '
Set-TextValue 'B11' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C11' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D11' '2.638.47'
Set-TextValue 'E11' '  +3.27%  '
Set-TextValue 'B12' 'TRON'
Set-TextValue 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D12' '0.104'
Set-TextValue 'E12' '  +0.77%  '
Set-TextValue 'D13' '15.82'
Set-TextValue 'E13' '  +1.03%  '
Set-TextValue 'D14' '23.87'
Set-TextValue 'E14' '  +9.89%  '
Set-TextValue 'E15' '  +3.30%  '
Set-TextValue 'E16' '  +1.88%  '
Set-TextValue 'D17' '2.303.29'
Set-TextValue 'E17' '  +3.57%  '
Set-TextValue 'D18' '43.809.53'
Set-TextValue 'E18' '  +5.23%  '
Set-TextValue 'D19' '0.0₃0930'
Set-TextValue 'E19' '  +4.49%  '
Set-TextValue 'D20' '73.37'
Set-TextValue 'E20' '  +0.89%  '
Set-TextValue 'D21' '6.25'
Set-TextValue 'E21' '  +3.74%  '
Set-TextValue 'D22' '250.10'
Set-TextValue 'E22' '  +0.25%  '
Set-TextValue 'E23' '  +0.17%  '
Set-TextValue 'E24' '  +7.15%  '
Set-TextValue 'E25' '  +2.27%  '
Set-TextValue 'D26' '9.87'
Set-TextValue 'E26' '  +3.42%  '
Set-TextValue 'D27' '169.93'
Set-TextValue 'E27' '  +1.38%  '
Set-TextValue 'E28' '  +0.39%  '
Set-TextValue 'D29' '20.60'
Set-TextValue 'E29' '  +3.46%  '
Set-TextValue 'D30' '1.48'
Set-TextValue 'E30' '  +5.44%  '
Set-TextValue 'E31' '  +0.99%  '
Set-TextValue 'E32' '  +0.17%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '4.75'
Set-TextValue 'E33' '  +2.82%  '
Set-TextValue 'B34' 'InternetComputer(DFINITY)'
Set-TextValue 'C34' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D34' '5.01'
Set-TextValue 'E34' '  +1.29%  '
Set-TextValue 'E35' '  +5.81%  '
Set-TextValue 'E36' '  +3.91%  '
Set-TextValue 'D37' '6.57'
Set-TextValue 'E37' '  -1.10%  '
Set-TextValue 'D38' '3.66'
Set-TextValue 'E38' '  -0.80%  '
Set-TextValue 'D39' '0.0251'
Set-TextValue 'E39' '  +4.51%  '
Set-TextValue 'E40' '  +0.28%  '
Set-TextValue 'D41' '8.80'
Set-TextValue 'E41' '  +1.93%  '
Set-TextValue 'D42' '4.59'
Set-TextValue 'E42' '  -4.28%  '
Set-TextValue 'D43' '0.0971'
Set-TextValue 'E43' '  -0.83%  '
Set-TextValue 'E44' '  -16.86%  '
Set-TextValue 'E45' '  +0.04%  '
Set-TextValue 'D46' '98.62'
Set-TextValue 'E46' '  +0.04%  '
Set-TextValue 'D47' '1.472.26'
Set-TextValue 'E47' '  +0.29%  '
Set-TextValue 'D48' '16.69'
Set-TextValue 'E48' '  +0.91%  '
Set-TextValue 'E49' '  +10.04%  '
Set-TextValue 'E50' '  +1.69%  '
Set-TextValue 'D51' '2.77'
Set-TextValue 'E51' '  -1.26%  '
